$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(86, '2026-02-15 02:55:40', '237671823369', 'MFS ENTREE COLLEGE MALANGUE', 7448),
    @(87, '2026-02-15 11:40:12', '237672128028', 'CAROLINE WAKO DJAMNOU', 34383),
    @(88, '2026-02-15 07:31:45', '237672277367', 'TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5', 28296),
    @(89, '2026-02-15 16:02:15', '237674853971', 'NJOSSEU TCHOUNZOU TOP MOBILE', 216629),
    @(90, '2026-02-15 17:24:40', '237674884705', 'BAH AMADOU MOUNTAGHA ETS MOBILE FINANCIAL SERVICES MFS', 13867),
    @(91, '2026-02-15 17:04:40', '237675779272', 'RODES NGWEM KEMAYOU', 19223),
    @(92, '2026-02-15 11:35:02', '237677304210', 'CARINE OROCK', 227157),
    @(93, '2026-02-06 09:50:11', '237678267353', 'LA NEGRESSE SARL EMBOLA BELTUS MBU', 0),
    @(94, '2026-02-15 14:41:54', '237678370615', 'ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS', 364192),
    @(95, '2026-02-15 13:39:21', '237678836319', 'KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS', 139993),
    @(96, '2026-02-15 14:28:30', '237678922502', 'NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT', 803074),
    @(97, '2026-02-15 16:16:41', '237679884264', 'MFS CICAM', 365895),
    @(98, '2026-02-15 12:10:35', '237681019523', 'ETS MOULAY RIPERT AND COMPANY', 401212),
    @(99, '2026-02-15 16:45:50', '237681125655', 'ANTOINETTE CLAUDE YOMBA', 95888),
    @(100, '2026-02-14 09:46:40', '237681240793', 'MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS', 2067),
    @(101, '2026-02-15 16:02:15', '237682117915', 'MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS', 287967),
    @(102, '2026-02-15 03:40:29', '237682154553', 'N A ISUFUH MIEMONA NGESSY ETS MOBILE FINANCIAL SERVICES MFS', 16500),
    @(103, '2026-02-15 13:09:15', '237682803277', 'NGUEMASSOM RENE MARTIAL LA NEGRESSE SARL', 67657),
    @(104, '2026-02-15 14:54:22', '237683323481', 'ETS LE CONTENT TSAZE DONFOUET FLORETTE ROSINE', 284338),
    @(105, '2026-02-15 14:40:50', '237683368985', 'RAISSA GOUEMKUM KENGNE', 39540),
    @(106, '2026-02-15 07:12:29', '237683432110', 'ERODINE TOUMENI', 552165),
    @(107, '2026-02-15 16:24:43', '237683743490', 'ETS LE CONTENT NGAH MARIE', 6753),
    @(108, '2026-02-15 11:10:49', '237651346686', 'DOUANLA ARMEL VISION TRADING COMPLEX AND TECHNOLOGIES SARL VISION TRADING COMPLEX', 13378),
    @(109, '2026-02-15 10:40:06', '237654508603', 'NINA MORELLE KUEDONG NGUETSA', 363989),
    @(110, '2026-02-15 13:12:27', '237672731004', 'DIALLO OUSMANE ETS MOBILE FINANCIAL SERVICES MFS', 1491),
    @(111, '2026-02-15 07:40:33', '237674884831', 'THOMAS NANFACK', 911655),
    @(112, '2026-02-15 15:46:35', '237675986754', 'VIVIANE MANTO DZOUA', 16304),
    @(113, '2026-02-15 06:11:00', '237678384388', 'DJOUMESSI KEUMOU FLAUBERT COCADIS SARL', 32976),
    @(114, '2026-02-15 09:25:07', '237678559161', 'ANCHANG LOUIS FUL ETS MOBILE FINANCIAL SERVICES MFS', 52443),
    @(115, '2026-02-15 15:47:14', '237678872943', 'ZANKIA RAISSA DIANE ETS MOBILE FINANCIAL SERVICES MFS', 91313),
    @(116, '2026-02-15 16:49:22', '237679244994', 'WEPNYU NJEBOH CLINTON CONTACT PLUS Sarl', 3019690),
    @(117, '2026-02-15 17:34:03', '237679508295', 'Valentine Kum Kang World T Plus', 123036),
    @(118, '2026-02-15 13:58:23', '237679808809', 'TONLEPEU NGALAGNI CHANCELINE DIANE ETS MOBILE FINANCIAL SERVICES MFS', 232441),
    @(119, '2026-02-15 15:27:27', '237679909537', 'LAURE GISELE TADOUMFO NOUOJI', 14867),
    @(120, '2026-02-15 12:11:04', '237681114370', 'ETS LE CONTENT SEDJINE TCHINDA RAMSES GAUTIER', 30),
    @(121, '2026-02-15 15:52:05', '237681490029', 'DOUNTIO FOFIE SAMANTA TOP MOBIL TELECOM', 91696),
    @(122, '2026-02-12 15:07:09', '237681655237', 'ETS LE CONTENT DJOUFACK WOUAFACK ALAIN GAROUSTE', 44713),
    @(123, '2026-02-15 16:53:40', '237681662606', 'GAEL PHALENNE NANA POUASSI', 108743),
    @(124, '2026-02-15 14:20:39', '237681662680', 'pascaline djiogo mejioguezem', 208573),
    @(125, '2026-02-15 15:49:04', '237681679214', 'N A NTOHMBU NKWIAMBOH FOUGUMA', 99703),
    @(126, '2026-02-15 16:42:24', '237651927707', 'JERRY MAFFO YPOLLI', 237591),
    @(127, '2026-02-15 16:43:00', '237654270800', 'MBANDA MARTIN BABIAKA ETS MOBILE FINANCIAL SERVICES MFS', 630199),
    @(128, '2026-02-14 14:09:07', '237674069453', 'DANIELLE SORELLE KENGNE  WENCI', 212893),
    @(129, '2026-02-15 15:21:35', '237674895078', 'ALLIANCE FLORE NOUMO EPSE TAGNE FOM', 50851),
    @(130, '2026-02-14 12:22:42', '237674895309', 'BORIS LIONEL FOTSING WENCI', 464939),
    @(131, '2026-02-15 09:33:20', '237675826419', 'Sara Gaetane Njouma epse Ntoma', 120628),
    @(132, '2026-02-15 11:24:23', '237676286294', 'NESTEN BOULANDI', 237800),
    @(133, '2026-02-15 13:26:30', '237677313421', 'NADEGE OMONINGO EPOUSE NEBA', 264809),
    @(134, '2026-02-15 13:55:16', '237677880357', 'vanissa sandjouon', 1209722),
    @(135, '2026-02-15 17:24:04', '237678201584', 'Fon Hortencia Engochan ONYEKA CONNECTION', 1162981),
    @(136, '2026-02-15 08:12:07', '237678239927', 'LA NEGRESSE SARL MADAGAL MOHAMED', 12423),
    @(137, '2026-02-15 11:32:58', '237678623874', 'BOUDIEU RHOTA KAMILAH CONNECTION', 328324),
    @(138, '2026-02-15 14:21:24', '237678637179', 'ALAIN GEORGES MOUDIO', 904324),
    @(139, '2026-02-15 06:45:11', '237679567513', 'DJOUTCHOU TCHANA MERVEILLE ESPERANCE', 17932)
)

$startRow = 86
$endRow = 139
$rangeAddr = "B" + $startRow + ":B" + $endRow
$ws.Range($rangeAddr).NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
